# Logboek: split the "sive gemaakt." run into three runs (with proofErr
# spell-check markers) and append the "Sprint 2" section with its first
# few log entries, exactly as described by the target diff.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 0: remove the existing "_GoBack" bookmark. It currently sits
# right after "sive gemaakt." (end of the document); we will re-add an
# equivalent bookmark at the new end of the document in step 2.
# ---------------------------------------------------------------------
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# ---------------------------------------------------------------------
# Step 1: within the last paragraph ("... About page gestyled en
# responsive gemaakt."), replace the trailing "sive gemaakt." text with
# three runs split around "gemaakt" (wrapped in spell-check proofErr
# markers), matching the target XML.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$searchRange = $d.Content
$searchRange.Start = $lastPara.Range.Start
$searchRange.End = $lastPara.Range.End
$null = $searchRange.Find.Execute("sive gemaakt.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Re-materialize a plain Range with the bounds Find located: InsertXML on a
# range that just had Find.Execute run on it directly tends to *insert
# after* instead of *replacing*, so we hand it a freshly constructed Range
# with the same Start/End instead.
$target = $d.Range($searchRange.Start, $searchRange.End)

$splitXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">sive </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>gemaakt</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($splitXml)

# ---------------------------------------------------------------------
# Step 2: append the new "Sprint 2" section after that paragraph: a
# blank line, the "Sprint 2" heading, the "31 mei" / "Gestart met CMS
# maken" entries, a blank line, and the "2 juni" / "CMS af en content
# uit db op site gezet (ipv. Zelf typen)" entries. The "_GoBack"
# bookmark is recreated at the very end, matching its original
# (unmoved, still-empty) position relative to the document content.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$insertionPoint = $d.Content
$insertionPoint.Start = $lastPara.Range.End
$insertionPoint.End = $lastPara.Range.End

$newSectionXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Kop1"/><w:rPr><w:sz w:val="40"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="40"/></w:rPr><w:t xml:space="preserve">Sprint </w:t></w:r><w:r><w:rPr><w:sz w:val="40"/></w:rPr><w:t>2</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">31 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>mei</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>Gestart</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> met CMS </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>maken</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">2 </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t>juni</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t xml:space="preserve">CMS af en content uit </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>d</w:t></w:r><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>b</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t xml:space="preserve"> op site gezet (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>ipv</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="32"/></w:rPr><w:t>. Zelf typen)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertionPoint.InsertXML($newSectionXml)
